$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Helper: convert a 1-based column index to its A1 column letters
# ------------------------------------------------------------------
function ColLetter([int]$n) {
    $s = ""
    while ($n -gt 0) {
        $rem = ($n - 1) % 26
        $s = [char](65 + $rem) + $s
        $n = [int](($n - $rem - 1) / 26)
    }
    return $s
}

# ====================================================================
# Sheet: Calculations
# ====================================================================
$calc = $wb.Worksheets.Item("Calculations")

# Column A width changed
$calc.Columns.Item(1).ColumnWidth = 17.9296875

# Row 5 label is unaffected in content (shared string re-numbering handles
# itself); nothing to change there explicitly.

# Row 6: years now run 2019 .. 2050 across B6:AG6 (was 2018..2050 across
# B6:AH6) -- shift the whole series left by one column and drop the old
# AH6 cell entirely.
for ($i = 0; $i -lt 32; $i++) {
    $calc.Cells.Item(6, 2 + $i).Value = 2019 + $i
}
$calc.Range("AH6").ClearContents()

# Row 7: new series label + new data (after storage & DR), one column
# shorter than before -- AH7 ends up blank (but keeps its style).
$calc.Range("A7").Value = "Peak Power Demand after Storage and DR[summer] : MostRecentRun"
$row7Values = @(752136,683006,728954,745991,757422,766081,772495,776739,779659,784465,792336,796481,800947,806133,810368,814654,819808,824754,830998,836502,841698,847887,853549,859603,866464,873306,881021,888875,896015,903925,912866,922440)
for ($i = 0; $i -lt $row7Values.Length; $i++) {
    $calc.Cells.Item(7, 2 + $i).Value = $row7Values[$i]
}
$calc.Range("AH7").ClearContents()

# Restore the view's selection on this sheet
$calc.Range("AH6:AH7").Select()

# ====================================================================
# Sheet: DRC-BDRC
# ====================================================================
$bdrc = $wb.Worksheets.Item("DRC-BDRC")

# D2:AH2 formulas shift their Calculations-row7 reference back one column
# and the fixed denominator moves from $C$7 to $B$7 (Calculations!B7 is
# now the 2019 figure).
for ($col = 4; $col -le 34; $col++) {
    $target = ColLetter $col
    $calcCol = ColLetter ($col - 1)
    $bdrc.Range("$target`2").Formula = "=`$C`$2*(Calculations!$calcCol`7/Calculations!`$B`$7)"
}

$bdrc.Range("D2:AH2").Select()

# ====================================================================
# Sheet: DRC-PADRC
# ====================================================================
$padrc = $wb.Worksheets.Item("DRC-PADRC")

# C2 now nets out the BAU portion already captured on DRC-BDRC
$padrc.Range("C2").Formula = "=Calculations!A3-'DRC-BDRC'!C2"

# D2:N2 (TREND-based) also net out the matching DRC-BDRC cell
for ($col = 4; $col -le 14; $col++) {
    $target = ColLetter $col
    $padrc.Range("$target`2").Formula = "=TREND(Calculations!`$A`$3:`$B`$3,Calculations!`$A`$2:`$B`$2,'DRC-PADRC'!$target`1)-'DRC-BDRC'!$target`2"
}

# O2:AH2 shift their Calculations-row7 reference back one column (same as
# DRC-BDRC), the fixed denominator moves from $N$7 to $M$7, and the
# matching DRC-BDRC cell is subtracted off.
for ($col = 15; $col -le 34; $col++) {
    $target = ColLetter $col
    $calcCol = ColLetter ($col - 1)
    $padrc.Range("$target`2").Formula = "=`$N`$2*(Calculations!$calcCol`7/Calculations!`$M`$7)-'DRC-BDRC'!$target`2"
}

$padrc.Range("O2:AH2").Select()

$wb.Application.Calculate()
Write-Output "edit applied"
